# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.195.50'
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').Value = '2.346.64'
$ws.Range('E3').Value = '  +6.06%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '''312.90'
$ws.Range('E5').Value = '  +6.04%  '
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = '''109.73'
$ws.Range('E6').Value = '  +1.56%  '
$ws.Range('D7').Value = '''0.641'
$ws.Range('E7').Value = '  +3.02%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = '''0.632'
$ws.Range('E9').Value = '  +6.21%  '
$ws.Range('D10').Value = '''42.81'
$ws.Range('E10').Value = '  -1.86%  '
$ws.Range('D11').Value = '''0.0938'
$ws.Range('E11').Value = '  +3.32%  '
$ws.Range('D12').Value = '''8.89'
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').Value = '''1.05'
$ws.Range('E13').Value = '  +9.99%  '
$ws.Range('D14').Value = '''0.105'
$ws.Range('E14').Value = '  +2.26%  '
$ws.Range('D15').Value = '''16.32'
$ws.Range('E15').Value = '  +9.61%  '
$ws.Range('D16').Value = '2.705.39'
$ws.Range('E16').Value = '  +6.29%  '
$ws.Range('D17').Value = '2.482.43'
$ws.Range('E17').Value = '  +11.37%  '
$ws.Range('D18').Value = '43.214.13'
$ws.Range('E18').Value = '  +2.17%  '
$ws.Range('E19').Value = '  +4.51%  '
$ws.Range('D20').Value = '''7.28'
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').Value = '''75.38'
$ws.Range('E21').Value = '  +4.09%  '
$ws.Range('E22').Value = '  +14.28%  '
$ws.Range('D23').Value = '''3.44'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '''253.15'
$ws.Range('E24').Value = '  +11.31%  '
$ws.Range('D25').Value = '''9.10'
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('D26').Value = '''12.02'
$ws.Range('E26').Value = '  +4.11%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '''39.43'
$ws.Range('E28').Value = '  +3.36%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.26'
$ws.Range('E29').Value = '  +1.61%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '''22.39'
$ws.Range('E30').Value = '  +7.51%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '''174.63'
$ws.Range('E31').Value = '  +0.70%  '
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('D33').Value = '''0.0933'
$ws.Range('E33').Value = '  +5.94%  '
$ws.Range('E34').Value = '  +8.72%  '
$ws.Range('E35').Value = '  +6.19%  '
$ws.Range('D36').Value = '''4.97'
$ws.Range('E36').Value = '  -1.56%  '
$ws.Range('D37').Value = '''0.0378'
$ws.Range('E37').Value = '  +4.82%  '
$ws.Range('D38').Value = '''4.14'
$ws.Range('E38').Value = '  -2.90%  '
$ws.Range('E39').Value = '  +1.62%  '
$ws.Range('D40').Value = '''2.69'
$ws.Range('E40').Value = '  +11.00%  '
$ws.Range('D41').Value = '''72.66'
$ws.Range('E41').Value = '  +2.72%  '
$ws.Range('E42').Value = '  +14.35%  '
$ws.Range('D43').Value = '''0.233'
$ws.Range('E43').Value = '  +1.68%  '
$ws.Range('D44').Value = '''12.86'
$ws.Range('E44').Value = '  +2.16%  '
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('D46').Value = '''5.63'
$ws.Range('E46').Value = '  +4.19%  '
$ws.Range('D47').Value = '''9.32'
$ws.Range('E47').Value = '  +10.27%  '
$ws.Range('D48').Value = '''110.81'
$ws.Range('E48').Value = '  +7.76%  '
$ws.Range('D49').Value = '''1.30'
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('E50').Value = '  +3.21%  '
$ws.Range('D51').Value = '''69.89'
$ws.Range('E51').Value = '  +5.19%  '
